# Update "想去人数" (column F) counts on several rows across all four
# sheets, reflecting a refreshed scrape of the source data
# (commit: "Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 265
$ws1.Range("F11").Value = 3108
$ws1.Range("F23").Value = 190
$ws1.Range("F26").Value = 211
$ws1.Range("F28").Value = 329
$ws1.Range("F31").Value = 126
$ws1.Range("F36").Value = 314
$ws1.Range("F37").Value = 1076
$ws1.Range("F38").Value = 5113
$ws1.Range("F39").Value = 564
$ws1.Range("F40").Value = 303
$ws1.Range("F41").Value = 168

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F17").Value = 5
$ws2.Range("F23").Value = 405
$ws2.Range("F26").Value = 741
$ws2.Range("F36").Value = 461
$ws2.Range("F44").Value = 26

# --- Sheet "本地生活" ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 659
$ws3.Range("F5").Value = 460

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 659
$ws4.Range("F3").Value = 265
$ws4.Range("F4").Value = 460
$ws4.Range("F14").Value = 3108
$ws4.Range("F23").Value = 5
$ws4.Range("F27").Value = 190
$ws4.Range("F30").Value = 211
$ws4.Range("F31").Value = 329
$ws4.Range("F36").Value = 405
$ws4.Range("F40").Value = 314
$ws4.Range("F41").Value = 1076
$ws4.Range("F42").Value = 5113
$ws4.Range("F44").Value = 564
$ws4.Range("F45").Value = 461
$ws4.Range("F46").Value = 303
$ws4.Range("F47").Value = 168
